$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the four taxonomy* header cells (row 1, columns C:F) to standardServiceType*
$ws.Range("C1").Value = "standardServiceTypeLetter"
$ws.Range("D1").Value = "standardServiceTypeParcel"
$ws.Range("E1").Value = "standardServiceTypeLetterKey"
$ws.Range("F1").Value = "standardServiceTypeParcelKey"

# Move the active cell / selection from G1 to F24
$ws.Range("F24").Select()

# Widen columns B-F (inputs chosen so the exported 1/6-char grid matches the
# target widths of 31.47 / 42.66 / 38.83 / 35.7 / 26.83); columns A, G, H stay at 18.33
$ws.Columns.Item(2).ColumnWidth = 30.6666666666667
$ws.Columns.Item(3).ColumnWidth = 41.8333333333333
$ws.Columns.Item(4).ColumnWidth = 38
$ws.Columns.Item(5).ColumnWidth = 34.8333333333333
$ws.Columns.Item(6).ColumnWidth = 26
